# Jets "Players Data" workbook update:
# C.Davis was removed from the roster, so his row is deleted from both the
# "Rushing" and "Receiving" sheets (the rows below shift up to fill the gap,
# and the now-unused "C.Davis" shared string is dropped on save).

$wb = $excel.ActiveWorkbook

# Rushing: C.Davis is row 11 (A11 = 9, B11 = "C.Davis").
$wsRushing = $wb.Worksheets.Item("Rushing")
$wsRushing.Rows(11).Delete()

# Receiving: C.Davis is row 6 (A6 = 4, B6 = "C.Davis").
$wsReceiving = $wb.Worksheets.Item("Receiving")
$wsReceiving.Rows(6).Delete()

# Final state: "Receiving" is the active sheet/tab, with cell selections left
# where editing ended on each sheet.
$wsRushing.Activate()
$wsRushing.Range("A11:F13").Select()

$wsReceiving.Activate()
$wsReceiving.Range("E13").Select()
